# Auto-generated Excel COM-interop script to apply leve-profit data updates
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 22689.777
$ws.Range("I19").Value = 100
$ws.Range("K19").Value = 100
$ws.Range("M19").Value = 75
# Row 135
$ws.Range("H135").Value = 37038504
$ws.Range("I135").Value = 40000784
$ws.Range("J135").Value = 10000
$ws.Range("K135").Value = 360007056
$ws.Range("L135").Value = 90000
$ws.Range("M135").Value = -360004521
$ws.Range("N135").Value = -95070

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4142.4287
$ws.Range("I61").Value = 3799.4
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3799.4
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3587.4
$ws.Range("N61").Value = -5424
# Row 74
$ws.Range("H74").Value = 2337.9412
$ws.Range("I74").Value = 2337.9412
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2337.9412
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1463.9412
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 2337.9412
$ws.Range("I77").Value = 2337.9412
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 11689.706
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -7321.706000000002
$ws.Range("N77").ClearContents()
# Row 80
$ws.Range("H80").Value = 42500
# Row 81
$ws.Range("H81").Value = 74500
$ws.Range("I81").Value = 70000
$ws.Range("J81").Value = 79000
$ws.Range("K81").Value = 70000
$ws.Range("L81").Value = 79000
$ws.Range("M81").Value = -69002
$ws.Range("N81").Value = -80996
# Row 83
$ws.Range("H83").Value = 42500
# Row 84
$ws.Range("H84").Value = 74500
$ws.Range("I84").Value = 70000
$ws.Range("J84").Value = 79000
$ws.Range("K84").Value = 210000
$ws.Range("L84").Value = 237000
$ws.Range("M84").Value = -205008
$ws.Range("N84").Value = -246984
# Row 122
$ws.Range("H122").Value = 9525913
$ws.Range("I122").Value = 13334865
$ws.Range("J122").Value = 3533.3
$ws.Range("K122").Value = 40004595
$ws.Range("L122").Value = 10599.9
$ws.Range("M122").Value = -40002145
$ws.Range("N122").Value = -15499.9
# Row 132
$ws.Range("H132").Value = 83336110
$ws.Range("I132").Value = 100002730
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 300008190
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -300005660
$ws.Range("N132").Value = -14060
# Row 136
$ws.Range("H136").Value = 4142.4287
$ws.Range("I136").Value = 3799.4
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11398.2
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -8848.200000000001
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1116.4324
$ws.Range("J105").Value = 1218.909
$ws.Range("L105").Value = 1218.909
$ws.Range("N105").Value = -4712.909
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 744339.4
$ws.Range("I132").Value = 530145.8
$ws.Range("J132").Value = 1253049
$ws.Range("K132").Value = 1590437.4
$ws.Range("L132").Value = 3759147
$ws.Range("M132").Value = -1587907.4
$ws.Range("N132").Value = -3764207
# Row 134
$ws.Range("H134").Value = 2853.5217
$ws.Range("I134").Value = 1983.5454
$ws.Range("J134").Value = 3651
$ws.Range("K134").Value = 5950.6362
$ws.Range("L134").Value = 10953
$ws.Range("M134").Value = -3415.6362
$ws.Range("N134").Value = -16023

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 9187.25
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 9666.4
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 28999.2
$ws.Range("M46").Value = -5909
$ws.Range("N46").Value = -29181.2
# Row 68
$ws.Range("H68").Value = 1912.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 1912.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 103
$ws.Range("H103").Value = 180.71428
$ws.Range("I103").Value = 160.83333
$ws.Range("J103").Value = 300
$ws.Range("K103").Value = 482.49999
$ws.Range("L103").Value = 900
$ws.Range("M103").Value = 396.50001
$ws.Range("N103").Value = -2658
# Row 106
$ws.Range("H106").Value = 4514.5
$ws.Range("J106").Value = 4514.5
$ws.Range("L106").Value = 13543.5
$ws.Range("N106").Value = -15435.5
# Row 121
$ws.Range("H121").Value = 66496.94500000001
$ws.Range("J121").Value = 89204.38
$ws.Range("L121").Value = 267613.14
$ws.Range("N121").Value = -270233.14

$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 15625
$ws.Range("J20").Value = 15625
$ws.Range("L20").Value = 15625
$ws.Range("N20").Value = -16115
# Row 122
$ws.Range("H122").Value = 11491365
$ws.Range("I122").Value = 151526.33
$ws.Range("J122").Value = 71430510
$ws.Range("K122").Value = 454578.99
$ws.Range("L122").Value = 214291530
$ws.Range("M122").Value = -452128.99
$ws.Range("N122").Value = -214296430
# Row 132
$ws.Range("H132").Value = 259942.05
$ws.Range("I132").Value = 315523.03
$ws.Range("J132").Value = 5857.5713
$ws.Range("K132").Value = 946569.0900000001
$ws.Range("L132").Value = 17572.7139
$ws.Range("M132").Value = -944039.0900000001
$ws.Range("N132").Value = -22632.7139

$ws = $wb.Worksheets.Item("LTW")
# Row 133
$ws.Range("H133").Value = 113199.664
$ws.Range("J133").Value = 113199.664
$ws.Range("L133").Value = 113199.664
$ws.Range("N133").Value = -118259.664
# Row 136
$ws.Range("H136").Value = 3724.5
$ws.Range("I136").Value = 1949.5
$ws.Range("J136").Value = 5499.5
$ws.Range("K136").Value = 5848.5
$ws.Range("L136").Value = 16498.5
$ws.Range("M136").Value = -3298.5
$ws.Range("N136").Value = -21598.5

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 2632.0466
$ws.Range("I136").Value = 1533.8276
$ws.Range("J136").Value = 4906.9287
$ws.Range("K136").Value = 4601.4828
$ws.Range("L136").Value = 14720.7861
$ws.Range("M136").Value = -2051.4828
$ws.Range("N136").Value = -19820.7861
# Row 140
$ws.Range("H140").Value = 84999.664
$ws.Range("J140").Value = 84999.664
$ws.Range("L140").Value = 84999.664
$ws.Range("N140").Value = -95359.664
